$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the formatting of the
# other header cells (bold, bordered, centered) by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF data columns for rows 2-10
$values = @(
  @(7, 7),
  @(8, 8),
  @(7, 8),
  @(5, 6),
  @(2, 3),
  @(5, 5),
  @(8, 8),
  @(6, 7),
  @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).Value = $values[$i][0]
  $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
